$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("D5").Value = 114675.7608660715
$ws.Range("E5").Value = -0.01641925629363952
$ws.Range("F5").Value = 0.2155139888333799
$ws.Range("G5").Value = -1.503494341713601
$ws.Range("H5").Value = 12.32542315777192

# Row 9
$ws.Range("D9").Value = 117954.3612058273
$ws.Range("E9").Value = -0.06251435975850676
$ws.Range("F9").Value = 0.3160540038138213
$ws.Range("G9").Value = -1.689019865453669
$ws.Range("H9").Value = 10.27623777756433

# Row 10
$ws.Range("D10").Value = 119501.1508225489
$ws.Range("E10").Value = -0.09805678537918844
$ws.Range("F10").Value = 0.3991945684870488
$ws.Range("G10").Value = -1.794603864231762
$ws.Range("H10").Value = 9.351016474581021

# Row 11
$ws.Range("D11").Value = 121386.6426153521
$ws.Range("E11").Value = -0.1683220395657837
$ws.Range("F11").Value = 0.6892054890913244
$ws.Range("G11").Value = -2.484535634006572
$ws.Range("H11").Value = 11.72790351028722

# Row 14
$ws.Range("D14").Value = 112492.0220184911
$ws.Range("E14").Value = -0.002954982496032888
$ws.Range("F14").Value = 0.141823251220435
$ws.Range("G14").Value = -0.5007239429912972
$ws.Range("H14").Value = 5.450011004942664
